$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.448.79'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.644.02'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '300.22'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3789'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.49'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3497'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08052'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.216'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.70%  '
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.09'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.284'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.247'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.62%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001210'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.642.55'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.26'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06980'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.626'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.39'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.41'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.462.41'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.412'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.994'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("E27").Value = '  -0.89%  '
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.180'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.59'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.840.29'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.860'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.139'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -4.69%  '
$ws.Range("E34").Value = '  -7.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9887'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02686'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -4.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08748'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.908'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2414'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06781'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.85'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.57%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6878'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.292'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.51'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6385'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -2.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.924'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.239'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '127.42'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07664'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.246'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.41%  '
